# Insert a new data row at row 694 (pushing existing rows 694:797 down to 695:798)
# and populate it with the new price-report record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("694:694").Insert()

$ws.Range("A694").Value = 10
$ws.Range("B694").Value = 'Vega Modelo de Temuco'
$ws.Range("C694").Value = 'La Araucanía'
$ws.Range("D694").Value = 45077
$ws.Range("E694").Value = 9
$ws.Range("F694").Value = 100112043
$ws.Range("G694").Value = 'Pepino ensalada'
$ws.Range("H694").Value = 'Sin especificar'
$ws.Range("I694").Value = 'Primera'
$ws.Range("J694").Value = 235
$ws.Range("K694").Value = 14000
$ws.Range("L694").Value = 15000
$ws.Range("M694").Value = 14468
$ws.Range("N694").Value = '$/caja 60 unidades'
$ws.Range("O694").Value = 'Región de Arica y Parinacota'
$ws.Range("P694").Value = 241
$ws.Range("Q694").Value = 60
$ws.Range("R694").Value = 'Hortaliza'
